$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain text, matching inlineStr source data,
# since many values look numeric (e.g. thousand-dot-grouped) and would
# otherwise be auto-converted to floating point numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.808.57'
$ws.Range('E2').Value = '  -4.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.628.88'
$ws.Range('E3').Value = '  -5.82%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.43'
$ws.Range('E5').Value = '  -2.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.64'
$ws.Range('E6').Value = '  +4.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.629.59'
$ws.Range('E7').Value = '  -5.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.622'
$ws.Range('E8').Value = '  -6.41%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.702'
$ws.Range('E10').Value = '  -5.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.158'
$ws.Range('E11').Value = '  -9.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.21'
$ws.Range('E12').Value = '  +4.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000287'
$ws.Range('E13').Value = '  -10.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.44'
$ws.Range('E14').Value = '  -7.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.224.76'
$ws.Range('E15').Value = '  -5.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.650.46'
$ws.Range('E16').Value = '  -5.18%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.126'
$ws.Range('E17').Value = '  -2.36%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.05'
$ws.Range('E18').Value = '  -9.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.60'
$ws.Range('E19').Value = '  -9.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.10'
$ws.Range('E20').Value = '  -7.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '67.644.06'
$ws.Range('E21').Value = '  -4.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '404.44'
$ws.Range('E22').Value = '  -7.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.48'
$ws.Range('E23').Value = '  -6.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.52'
$ws.Range('E24').Value = '  -6.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.98'
$ws.Range('E25').Value = '  -9.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.04'
$ws.Range('E26').Value = '  +0.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.53'
$ws.Range('E27').Value = '  -9.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.61'
$ws.Range('E28').Value = '  -8.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.02'
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.32'
$ws.Range('E30').Value = '  -11.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.23'
$ws.Range('E31').Value = '  -7.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.10'
$ws.Range('E32').Value = '  -14.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.18'
$ws.Range('E33').Value = '  -9.62%  '
$ws.Range('E34').Value = '  -7.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '64.11'
$ws.Range('E35').Value = '  -6.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '42.52'
$ws.Range('E36').Value = '  -11.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '590.92'
$ws.Range('E37').Value = '  -8.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0875'
$ws.Range('E38').Value = '  -11.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.392'
$ws.Range('E40').Value = '  -9.39%  '
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.134'
$ws.Range('E42').Value = '  -8.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.97'
$ws.Range('E43').Value = '  -8.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.66'
$ws.Range('E44').Value = '  -9.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0430'
$ws.Range('E45').Value = '  -7.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.81'
$ws.Range('E46').Value = '  -13.04%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.69'
$ws.Range('E47').Value = '  -3.91%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.133'
$ws.Range('E48').Value = '  -7.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.87'
$ws.Range('E49').Value = '  -10.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.10'
$ws.Range('E50').Value = '  -5.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.682.06'
$ws.Range('E51').Value = '  -7.84%  '
